$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tier 1_obs")

for ($row = 6; $row -le 34; $row++) {
    $ws.Cells.Item($row, 12).Value = "not applicable"
}
